$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C29").Value = "Controllers and views"
$ws.Range("D29").Value = 44742
$ws.Range("E29").Value = "Creacion de controller y views de user"

$ws.Range("D27").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E29:E30").Select()
